# Applies attendance updates to Sheet1 of the active workbook.
# For each listed row, set the specified column cells from 0 to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    3  = @("G", "H")
    4  = @("D", "E")
    5  = @("D", "E")
    6  = @("D", "E")
    7  = @("H")
    8  = @("H")
    9  = @("D", "E")
    10 = @("D", "E")
    11 = @("D", "E")
    12 = @("H")
    13 = @("H")
    14 = @("H")
    15 = @("H")
    16 = @("H")
    17 = @("H")
    18 = @("H")
}

foreach ($row in $changes.Keys) {
    foreach ($col in $changes[$row]) {
        $ws.Range("$col$row").Value = 1
    }
}
